# Apply the "add priority index" + new supervisor rows update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows right after the header row (row 1) so the existing
# supervisor rows (Prof Catie ... Prof Harry) shift from rows 2-7 down to rows 6-11.
$ws.Range("A2:A5").EntireRow.Insert()

# New header for the priority column.
$ws.Range("D1").Value = "priority"

# First new supervisor row (row 2), filled across.
$ws.Range("A2").Value = "Prof Chan"
$ws.Range("B2").Value = "tid00001"
$ws.Range("C2").Value = "tpw00001"

# Remaining new supervisor rows (3-5), filled column by column.
$ws.Range("A3").Value = "Prof Lam"
$ws.Range("A4").Value = "Dr Amantha"
$ws.Range("A5").Value = "Dr Banana"

$ws.Range("B3").Value = "tid00002"
$ws.Range("B4").Value = "tid00003"
$ws.Range("B5").Value = "tid00004"

$ws.Range("C3").Value = "tpw00002"
$ws.Range("C4").Value = "tpw00003"
$ws.Range("C5").Value = "tpw00004"

# Priority values for the new supervisors.
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 3

# Priority values for the pre-existing supervisors, now shifted to rows 6-11.
$ws.Range("D6").Value = 2
$ws.Range("D7").Value = 3
$ws.Range("D8").Value = 3
$ws.Range("D9").Value = 2
$ws.Range("D10").Value = 3
$ws.Range("D11").Value = 1

$ws.Range("D5").Select()
$excel.ActiveWindow.Zoom = 139
